# "Actualizar" automation re-run: refresh the existing 44230.87725172282
# timestamps (rows 16-29) to 44230.87725172454, then append a brand-new
# availability snapshot (rows 30-43) mirroring the same Name/URL/Status
# pattern as rows 16-29, stamped with the newer check time
# 44230.88231832047, including matching hyperlinks on column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the timestamp already stored in D16:D29 -------------------
for ($r = 16; $r -le 29; $r++) {
    $ws.Range("D$r").Value = 44230.87725172454
}

# --- 2. Append the new block of rows 30-43 ---------------------------------
$names = @(
    "Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", "UtilidadesOdoo",
    "Filtros Dashboard", "MapStore", "GeoServer", "Tomcat", "Shiny", "Github", "EZ Exporter"
)

# Text written into column B (what the user sees / what is stored as the cell value)
$urlLabels = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Hyperlink target (Address) - the MapStore entry's visible "#/" is a sub-address/location
$hyperlinkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$hyperlinkLocations = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")

$newTimestamp = 44230.88231832047

for ($i = 0; $i -lt 14; $i++) {
    $row = 30 + $i

    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("C$row").Value = "Disponible"
    $ws.Range("D$row").Value = $newTimestamp
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($hyperlinkLocations[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $hyperlinkAddresses[$i], $hyperlinkLocations[$i])
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $hyperlinkAddresses[$i])
    }

    # Hyperlinks.Add() stamps the display text itself; overwrite with the
    # exact label text and re-apply the workbook's existing Hyperlink style
    # (same look as the other rows, col B).
    $ws.Range("B$row").Value = $urlLabels[$i]
    $ws.Range("B$row").Style = $ws.Range("B2").Style
}
